$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for the cells we touch so numeric-looking strings
# (e.g. "10" -> "11", "243.68" -> "243.32") keep their exact text representation
# instead of being auto-converted to numbers by Excel.
$targetRanges = @("D2","G2","D3","G3","D4","G4","D5","G5","D6","G6","D7","G7","G8","D9","G9","D10","G10","D11","G11","D12","G12","B13","C13","D13","E13","G13","B14","C14","D14","E14","G14","B15","C15","D15","E15","G15","B16","C16","D16","E16","G16","B17","C17","D17","E17","G17","B18","C18","D18","E18","G18","B19","C19","D19","E19","G19","B20","C20","D20","E20","G20","B21","C21","D21","E21","G21","B22","C22","D22","E22","G22","B23","C23","D23","E23","G23","B24","C24","D24","E24","G24","B25","C25","D25","E25","G25","B26","C26","D26","E26","G26","D27","G27","G28","G29","G30","G31","G32","G33","G34","G35","G36","G37","G38","G39","D40","G40","D41","G41","D42","G42","D43","E43","G43","D44","G44","D45","G45","G46","D47","E47","G47","D48","G48","D49","G49","D50","G50","G51")
foreach ($ref in $targetRanges) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value = "243.32"
$ws.Range("G2").Value = "11"
$ws.Range("D3").Value = "23.02"
$ws.Range("G3").Value = "11"
$ws.Range("D4").Value = "5.399"
$ws.Range("G4").Value = "11"
$ws.Range("D5").Value = "0.05916"
$ws.Range("G5").Value = "11"
$ws.Range("D6").Value = "3.459"
$ws.Range("G6").Value = "11"
$ws.Range("D7").Value = "6.541"
$ws.Range("G7").Value = "11"
$ws.Range("G8").Value = "11"
$ws.Range("D9").Value = "0.9105"
$ws.Range("G9").Value = "11"
$ws.Range("D10").Value = "0.1414"
$ws.Range("G10").Value = "11"
$ws.Range("D11").Value = "0.07372"
$ws.Range("G11").Value = "11"
$ws.Range("D12").Value = "0.03276"
$ws.Range("G12").Value = "11"
$ws.Range("B13").Value = "ProBitToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D13").Value = "0.1324"
$ws.Range("E13").Value = "12ProBitTokenPROB"
$ws.Range("G13").Value = "11"
$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D14").Value = "0.03074"
$ws.Range("E14").Value = "13BitrueCoinBTR"
$ws.Range("G14").Value = "11"
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").Value = "0.09355"
$ws.Range("E15").Value = "14BitMartTokenBMX"
$ws.Range("G15").Value = "11"
$ws.Range("B16").Value = "MCDex"
$ws.Range("C16").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D16").Value = "3.855"
$ws.Range("E16").Value = "15MCDexMCB"
$ws.Range("G16").Value = "11"
$ws.Range("B17").Value = "BitForexToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D17").Value = "0.001564"
$ws.Range("E17").Value = "16BitForexTokenBF"
$ws.Range("G17").Value = "11"
$ws.Range("B18").Value = "CoinExToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D18").Value = "0.04683"
$ws.Range("E18").Value = "17CoinExTokenCET"
$ws.Range("G18").Value = "11"
$ws.Range("B19").Value = "One"
$ws.Range("C19").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D19").Value = "0.0005943"
$ws.Range("E19").Value = "18OneONE"
$ws.Range("G19").Value = "11"
$ws.Range("B20").Value = "TigerCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D20").Value = "0.006035"
$ws.Range("E20").Value = "19TigerCashTCH"
$ws.Range("G20").Value = "11"
$ws.Range("B21").Value = "HotbitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D21").Value = "0.004979"
$ws.Range("E21").Value = "20HotbitTokenHTBBestin24h"
$ws.Range("G21").Value = "11"
$ws.Range("B22").Value = "BitKan"
$ws.Range("C22").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D22").Value = "0.0009821"
$ws.Range("E22").Value = "21BitKanKAN"
$ws.Range("G22").Value = "11"
$ws.Range("B23").Value = "NitroEx"
$ws.Range("C23").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D23").Value = "0.00008604"
$ws.Range("E23").Value = "22NitroExNTX"
$ws.Range("G23").Value = "11"
$ws.Range("B24").Value = "LEO"
$ws.Range("C24").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D24").Value = "3.610"
$ws.Range("E24").Value = "23LEOLEO"
$ws.Range("G24").Value = "11"
$ws.Range("B25").Value = "BTSEToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D25").Value = "2.151"
$ws.Range("E25").Value = "24BTSETokenBTSE"
$ws.Range("G25").Value = "11"
$ws.Range("B26").Value = "BitpandaEcosystemToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D26").Value = "0.3240"
$ws.Range("E26").Value = "25BitpandaEcosystemTokenBEST"
$ws.Range("G26").Value = "11"
$ws.Range("D27").Value = "0.0002902"
$ws.Range("G27").Value = "11"
$ws.Range("G28").Value = "11"
$ws.Range("G29").Value = "11"
$ws.Range("G30").Value = "11"
$ws.Range("G31").Value = "11"
$ws.Range("G32").Value = "11"
$ws.Range("G33").Value = "11"
$ws.Range("G34").Value = "11"
$ws.Range("G35").Value = "11"
$ws.Range("G36").Value = "11"
$ws.Range("G37").Value = "11"
$ws.Range("G38").Value = "11"
$ws.Range("G39").Value = "11"
$ws.Range("D40").Value = "0.03962"
$ws.Range("G40").Value = "11"
$ws.Range("D41").Value = "0.006193"
$ws.Range("G41").Value = "11"
$ws.Range("D42").Value = "0.1075"
$ws.Range("G42").Value = "11"
$ws.Range("D43").Value = "0.002621"
$ws.Range("E43").Value = "42CEJICEJIWorstin24h"
$ws.Range("G43").Value = "11"
$ws.Range("D44").Value = "0.008610"
$ws.Range("G44").Value = "11"
$ws.Range("D45").Value = "0.00005178"
$ws.Range("G45").Value = "11"
$ws.Range("G46").Value = "11"
$ws.Range("D47").Value = "0.8095"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOIN"
$ws.Range("G47").Value = "11"
$ws.Range("D48").Value = "0.002335"
$ws.Range("G48").Value = "11"
$ws.Range("D49").Value = "0.00002101"
$ws.Range("G49").Value = "11"
$ws.Range("D50").Value = "0.0002001"
$ws.Range("G50").Value = "11"
$ws.Range("G51").Value = "11"
